$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(51, 8).Value = 4522.933
$ws.Cells.Item(51, 10).Value = 4603.3076
$ws.Cells.Item(51, 12).Value = 4603.3076
$ws.Cells.Item(51, 14).Value = -5571.3076
$ws.Cells.Item(112, 8).Value = 2618.3635
$ws.Cells.Item(112, 10).Value = 3265.625
$ws.Cells.Item(112, 12).Value = 9796.875
$ws.Cells.Item(112, 14).Value = -12012.875
$ws.Cells.Item(132, 8).Value = 1169.5927
$ws.Cells.Item(132, 9).Value = 1107.7273
$ws.Cells.Item(132, 11).Value = 3323.1819
$ws.Cells.Item(132, 13).Value = -793.1819
$ws.Cells.Item(137, 8).Value = 2233.7334
$ws.Cells.Item(137, 9).Value = 2337.3
$ws.Cells.Item(137, 10).Value = 2026.6
$ws.Cells.Item(137, 11).Value = 7011.900000000001
$ws.Cells.Item(137, 12).Value = 6079.799999999999
$ws.Cells.Item(137, 13).Value = -4461.900000000001
$ws.Cells.Item(137, 14).Value = -11179.8

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2810.875
$ws.Cells.Item(2, 9).Value = 1720.4286
$ws.Cells.Item(2, 10).Value = 10444
$ws.Cells.Item(2, 11).Value = 1720.4286
$ws.Cells.Item(2, 12).Value = 10444
$ws.Cells.Item(2, 13).Value = -1607.4286
$ws.Cells.Item(2, 14).Value = -10670
$ws.Cells.Item(32, 8).Value = 3541.8867
$ws.Cells.Item(32, 9).Value = 2147.0977
$ws.Cells.Item(32, 11).Value = 2147.0977
$ws.Cells.Item(32, 13).Value = -1860.0977
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 14).ClearContents()
$ws.Cells.Item(61, 8).Value = 1463.1
$ws.Cells.Item(61, 9).Value = 1070.1111
$ws.Cells.Item(61, 11).Value = 1070.1111
$ws.Cells.Item(61, 13).Value = -858.1111000000001
$ws.Cells.Item(63, 8).Value = 3481.4
$ws.Cells.Item(63, 9).Value = 3799
$ws.Cells.Item(63, 10).Value = 3005
$ws.Cells.Item(63, 11).Value = 3799
$ws.Cells.Item(63, 12).Value = 3005
$ws.Cells.Item(63, 13).Value = -3113
$ws.Cells.Item(63, 14).Value = -4377
$ws.Cells.Item(66, 8).Value = 3481.4
$ws.Cells.Item(66, 9).Value = 3799
$ws.Cells.Item(66, 10).Value = 3005
$ws.Cells.Item(66, 11).Value = 18995
$ws.Cells.Item(66, 12).Value = 15025
$ws.Cells.Item(66, 13).Value = -15563
$ws.Cells.Item(66, 14).Value = -21889
$ws.Cells.Item(74, 8).Value = 3093.9143
$ws.Cells.Item(74, 9).Value = 3055.9092
$ws.Cells.Item(74, 11).Value = 3055.9092
$ws.Cells.Item(74, 13).Value = -2181.9092
$ws.Cells.Item(77, 8).Value = 3093.9143
$ws.Cells.Item(77, 9).Value = 3055.9092
$ws.Cells.Item(77, 11).Value = 15279.546
$ws.Cells.Item(77, 13).Value = -10911.546
$ws.Cells.Item(116, 8).Value = 2810.875
$ws.Cells.Item(116, 9).Value = 1720.4286
$ws.Cells.Item(116, 10).Value = 10444
$ws.Cells.Item(116, 11).Value = 1720.4286
$ws.Cells.Item(116, 12).Value = 10444
$ws.Cells.Item(116, 13).Value = 573.5714
$ws.Cells.Item(116, 14).Value = -15032
$ws.Cells.Item(136, 8).Value = 1463.1
$ws.Cells.Item(136, 9).Value = 1070.1111
$ws.Cells.Item(136, 11).Value = 3210.3333
$ws.Cells.Item(136, 13).Value = -660.3333000000002

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2810.875
$ws.Cells.Item(3, 9).Value = 1720.4286
$ws.Cells.Item(3, 10).Value = 10444
$ws.Cells.Item(3, 11).Value = 1720.4286
$ws.Cells.Item(3, 12).Value = 10444
$ws.Cells.Item(3, 13).Value = -1606.4286
$ws.Cells.Item(3, 14).Value = -10672
$ws.Cells.Item(134, 8).Value = 3181.4917
$ws.Cells.Item(134, 9).Value = 3009
$ws.Cells.Item(134, 11).Value = 9027
$ws.Cells.Item(134, 13).Value = -6492

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 4239.4546
$ws.Cells.Item(16, 9).Value = 6088.6
$ws.Cells.Item(16, 11).Value = 6088.6
$ws.Cells.Item(16, 13).Value = -5801.6
$ws.Cells.Item(99, 8).Value = 4391.091
$ws.Cells.Item(99, 9).Value = 2503.6667
$ws.Cells.Item(99, 11).Value = 2503.6667
$ws.Cells.Item(99, 13).Value = -1005.6667
$ws.Cells.Item(113, 8).Value = 4239.4546
$ws.Cells.Item(113, 9).Value = 6088.6
$ws.Cells.Item(113, 11).Value = 6088.6
$ws.Cells.Item(113, 13).Value = -3918.6
$ws.Cells.Item(126, 8).Value = 4391.091
$ws.Cells.Item(126, 9).Value = 2503.6667
$ws.Cells.Item(126, 11).Value = 7511.000100000001
$ws.Cells.Item(126, 13).Value = -5041.000100000001
$ws.Cells.Item(132, 8).Value = 2465.5386
$ws.Cells.Item(132, 9).Value = 1410.9333
$ws.Cells.Item(132, 11).Value = 4232.7999
$ws.Cells.Item(132, 13).Value = -1702.7999

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 13000
$ws.Cells.Item(5, 9).Value = 13000
$ws.Cells.Item(5, 11).Value = 13000
$ws.Cells.Item(5, 13).Value = -12888
$ws.Cells.Item(6, 8).Value = 0
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 12).Value = 0
$ws.Cells.Item(6, 14).Value = 0
$ws.Cells.Item(9, 8).Value = 2807.6667
$ws.Cells.Item(9, 9).Value = 1655
$ws.Cells.Item(9, 10).Value = 4248.5
$ws.Cells.Item(9, 11).Value = 1655
$ws.Cells.Item(9, 12).Value = 4248.5
$ws.Cells.Item(9, 13).Value = -1485
$ws.Cells.Item(9, 14).Value = -4588.5
$ws.Cells.Item(13, 8).Value = 547.8333
$ws.Cells.Item(13, 9).Value = 322.5
$ws.Cells.Item(13, 10).Value = 998.5
$ws.Cells.Item(13, 11).Value = 322.5
$ws.Cells.Item(13, 12).Value = 998.5
$ws.Cells.Item(13, 13).Value = -183.5
$ws.Cells.Item(13, 14).Value = -1276.5
$ws.Cells.Item(14, 8).Value = 2363.75
$ws.Cells.Item(14, 9).Value = 2225
$ws.Cells.Item(14, 10).Value = 2502.5
$ws.Cells.Item(14, 11).Value = 2225
$ws.Cells.Item(14, 12).Value = 2502.5
$ws.Cells.Item(14, 13).Value = -2057
$ws.Cells.Item(14, 14).Value = -2838.5
$ws.Cells.Item(16, 8).Value = 0
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 14).Value = 0
$ws.Cells.Item(17, 8).Value = 6200
$ws.Cells.Item(17, 10).Value = 6200
$ws.Cells.Item(17, 12).Value = 6200
$ws.Cells.Item(17, 14).Value = -6536
$ws.Cells.Item(22, 8).Value = 3272.7273
$ws.Cells.Item(22, 9).Value = 4751
$ws.Cells.Item(22, 10).Value = 2428
$ws.Cells.Item(22, 11).Value = 4751
$ws.Cells.Item(22, 12).Value = 2428
$ws.Cells.Item(22, 13).Value = -4222
$ws.Cells.Item(22, 14).Value = -3486
$ws.Cells.Item(23, 8).Value = 0
$ws.Cells.Item(23, 10).Value = 0
$ws.Cells.Item(23, 12).Value = 0
$ws.Cells.Item(23, 14).Value = 0
$ws.Cells.Item(36, 8).Value = 8000
$ws.Cells.Item(36, 10).Value = 0
$ws.Cells.Item(36, 12).Value = 0
$ws.Cells.Item(36, 14).ClearContents()

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 44111.332
$ws.Cells.Item(22, 9).Value = 167362.5
$ws.Cells.Item(22, 10).Value = 3027.611
$ws.Cells.Item(22, 11).Value = 167362.5
$ws.Cells.Item(22, 12).Value = 3027.611
$ws.Cells.Item(22, 13).Value = -167067.5
$ws.Cells.Item(22, 14).Value = -3617.611
$ws.Cells.Item(27, 8).Value = 44111.332
$ws.Cells.Item(27, 9).Value = 167362.5
$ws.Cells.Item(27, 10).Value = 3027.611
$ws.Cells.Item(27, 11).Value = 167362.5
$ws.Cells.Item(27, 12).Value = 3027.611
$ws.Cells.Item(27, 13).Value = -167255.5
$ws.Cells.Item(27, 14).Value = -3241.611
$ws.Cells.Item(55, 8).Value = 657.913
$ws.Cells.Item(55, 10).Value = 1511
$ws.Cells.Item(55, 12).Value = 1511
$ws.Cells.Item(55, 14).Value = -1857
$ws.Cells.Item(61, 8).Value = 1661.4
$ws.Cells.Item(61, 10).Value = 5005
$ws.Cells.Item(61, 12).Value = 5005
$ws.Cells.Item(61, 14).Value = -5409
$ws.Cells.Item(113, 8).Value = 1661.4
$ws.Cells.Item(113, 10).Value = 5005
$ws.Cells.Item(113, 12).Value = 5005
$ws.Cells.Item(113, 14).Value = -9345
$ws.Cells.Item(132, 8).Value = 3810.4736
$ws.Cells.Item(132, 10).Value = 1960.6666
$ws.Cells.Item(132, 12).Value = 5881.9998
$ws.Cells.Item(132, 14).Value = -10941.9998
$ws.Cells.Item(136, 8).Value = 2456.5366
$ws.Cells.Item(136, 9).Value = 2061.75
$ws.Cells.Item(136, 10).Value = 5299
$ws.Cells.Item(136, 11).Value = 6185.25
$ws.Cells.Item(136, 12).Value = 15897
$ws.Cells.Item(136, 13).Value = -3635.25
$ws.Cells.Item(136, 14).Value = -20997
